$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data per the latest scrape.
# Numeric-looking "Price" strings are prefixed with a literal leading
# apostrophe so Excel stores them as text (matching the original
# t="inlineStr" / shared-string cell type) instead of coercing them to numbers.

$ws.Range("D2").Value = '58.224.90'
$ws.Range("E2").Value = '  +1.18%  '
$ws.Range("D3").Value = '2.353.53'
$ws.Range("E3").Value = '  +1.53%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Value = '''540.25'
$ws.Range("E5").Value = '  +1.58%  '
$ws.Range("D6").Value = '''135.86'
$ws.Range("E6").Value = '  +2.43%  '
$ws.Range("E8").Value = '  +4.92%  '
$ws.Range("E9").Value = '  +0.24%  '
$ws.Range("E10").Value = '  +5.37%  '
$ws.Range("E11").Value = '  -0.84%  '
$ws.Range("E12").Value = '  +2.35%  '
$ws.Range("D13").Value = '''23.84'
$ws.Range("E13").Value = '  +1.53%  '
$ws.Range("D14").Value = '2.772.26'
$ws.Range("E14").Value = '  +1.27%  '
$ws.Range("D15").Value = '58.186.69'
$ws.Range("E15").Value = '  +1.34%  '
$ws.Range("E16").Value = '  +0.46%  '
$ws.Range("D17").Value = '2.332.54'
$ws.Range("E17").Value = '  -0.30%  '
$ws.Range("E18").Value = '  +2.88%  '
$ws.Range("D19").Value = '''332.28'
$ws.Range("E19").Value = '  -1.73%  '
$ws.Range("E20").Value = '  +2.66%  '
$ws.Range("E21").Value = '  -1.58%  '
$ws.Range("E22").Value = '  +0.06%  '
$ws.Range("D23").Value = '''62.84'
$ws.Range("E23").Value = '  +1.13%  '
$ws.Range("E24").Value = '  -0.07%  '
$ws.Range("B25").Value = 'InternetComputer(DFINITY)'
$ws.Range("C25").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D25").Value = '''8.49'
$ws.Range("E25").Value = '  -2.52%  '
$ws.Range("B26").Value = 'Binance-PegBSC-USD'
$ws.Range("C26").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D26").Value = '''1.00'
$ws.Range("E26").Value = '  +0.70%  '
$ws.Range("E27").Value = '  +3.49%  '
$ws.Range("D28").Value = '''172.11'
$ws.Range("E28").Value = '  -0.65%  '
$ws.Range("E29").Value = '  +1.41%  '
$ws.Range("E30").Value = '  +1.62%  '
$ws.Range("E31").Value = '  +0.25%  '
$ws.Range("E32").Value = '  +12.40%  '
$ws.Range("D33").Value = '''18.43'
$ws.Range("E33").Value = '  -0.63%  '
$ws.Range("E34").Value = '  +0.03%  '
$ws.Range("D35").Value = '''4.27'
$ws.Range("E35").Value = '  +6.82%  '
$ws.Range("D36").Value = '''1.00'
$ws.Range("E36").Value = '  +0.25%  '
$ws.Range("E37").Value = '  -0.54%  '
$ws.Range("D38").Value = '''1.65'
$ws.Range("E38").Value = '  +4.00%  '
$ws.Range("D39").Value = '''39.20'
$ws.Range("E39").Value = '  -0.16%  '
$ws.Range("D40").Value = '''145.13'
$ws.Range("E40").Value = '  -2.58%  '
$ws.Range("D41").Value = '''293.58'
$ws.Range("E41").Value = '  +4.53%  '
$ws.Range("D42").Value = '''0.378'
$ws.Range("E42").Value = '  +0.87%  '
$ws.Range("D43").Value = '''3.65'
$ws.Range("E43").Value = '  +1.14%  '
$ws.Range("E44").Value = '  +1.86%  '
$ws.Range("D45").Value = '''19.22'
$ws.Range("E45").Value = '  +1.71%  '
$ws.Range("D46").Value = '''0.0503'
$ws.Range("E46").Value = '  +0.26%  '
$ws.Range("D47").Value = '''0.563'
$ws.Range("E47").Value = '  +0.82%  '
$ws.Range("E48").Value = '  +1.43%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '''17.53'
$ws.Range("E49").Value = '  +0.39%  '
$ws.Range("B50").Value = 'Polygon'
$ws.Range("C50").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D50").Value = '''0.382'
$ws.Range("E50").Value = '  +0.01%  '
$ws.Range("E51").Value = '  +0.39%  '
